# Insert 3 new data rows (Especial / Primera / Segunda) right before the
# current row 569, shifting the existing rows 569:672 down to 572:675.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("569:571").Insert()

$values = @(
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44641, 4, "Fruta", 100108, "Tropicales y subtropicales", 100108002, "Mango", "Sin especificar", "Especial", 512, 7500, 8000, 7750, "`$/bandeja 4 kilos", "Perú", 1938, 4),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44641, 4, "Fruta", 100108, "Tropicales y subtropicales", 100108002, "Mango", "Sin especificar", "Primera", 512, 7500, 8000, 7750, "`$/bandeja 4 kilos", "Perú", 1938, 4),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44641, 4, "Fruta", 100108, "Tropicales y subtropicales", 100108002, "Mango", "Sin especificar", "Segunda", 512, 7500, 8000, 7750, "`$/bandeja 4 kilos", "Perú", 1938, 4)
)

for ($i = 0; $i -lt 3; $i++) {
    $r = 569 + $i
    $row = $values[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
    $ws.Cells.Item($r, 14).Value = $row[13]
    $ws.Cells.Item($r, 15).Value = $row[14]
    $ws.Cells.Item($r, 16).Value = $row[15]
    $ws.Cells.Item($r, 17).Value = $row[16]
    $ws.Cells.Item($r, 18).Value = $row[17]
    $ws.Cells.Item($r, 19).Value = $row[18]
    $ws.Cells.Item($r, 20).Value = $row[19]
}
